$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.695.53'
$ws.Range('E2').Value = '  -1.27%  '

$ws.Range('D3').Value = '2.029.15'
$ws.Range('E3').Value = '  -1.39%  '

$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.32'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.15%  '

$ws.Range('E6').Value = '  -1.38%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.08'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.41%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('E9').Value = '  -2.46%  '

$ws.Range('E10').Value = '  +1.34%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.104'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.47%  '

$ws.Range('D12').Value = '2.331.18'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.41'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.26%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.01'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.38%  '

$ws.Range('E15').Value = '  +0.04%  '

$ws.Range('E16').Value = '  -2.43%  '

$ws.Range('D17').Value = '2.040.13'
$ws.Range('E17').Value = '  -0.61%  '

$ws.Range('D18').Value = '37.650.09'
$ws.Range('E18').Value = '  -1.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.47'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.75%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.91'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -6.06%  '

$ws.Range('D21').Value = '0.0₃0820'
$ws.Range('E21').Value = '  -1.67%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.21'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.13%  '

$ws.Range('E23').Value = '  +0.58%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.41'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.99%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.24'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.24%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.47'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.58%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.22'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.77%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.129'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.97%  '

$ws.Range('E29').Value = '  -1.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.26'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.67%  '

$ws.Range('E31').Value = '  +0.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.21'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.45%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.38'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.67%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0602'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.66%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.46'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.94%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.32'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.63%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.28'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.26%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.34'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.67%  '

$ws.Range('E39').Value = '  +0.06%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.74'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.24%  '

$ws.Range('D41').Value = '1.540.90'
$ws.Range('E41').Value = '  +1.04%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0215'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '95.48'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.62%  '

$ws.Range('E44').Value = '  -2.57%  '

$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.25'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.96%  '

$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0910'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.26%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.10'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.92%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.85%  '

$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.96'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.52%  '

$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.10'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.44%  '

$ws.Range('D51').Value = '2.221.10'
$ws.Range('E51').Value = '  -1.25%  '
